# Add new serotypes "36A" and "36B" to the "1. Serotype" sheet, replacing the
# old single "36" row (row 66), which pushes every following row down by one.
#
# Also updates the workbook view state so that the "1. Serotype" sheet is the
# active/selected tab (scrolled down near the newly-inserted rows), instead
# of "4. Variants".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("1. Serotype")

# --- Data edit -----------------------------------------------------------
# Row 66 currently holds the single serotype "36". Duplicate that row
# (copy + insert keeps the cell styles/types, e.g. the "False" text in
# column C stays a shared string instead of turning into a boolean), then
# relabel the original row as "36A" and the newly inserted row as "36B".
# Every row below (37, 38, 39, ...) shifts down by one automatically.
$ws1.Rows.Item(66).Copy()
$ws1.Rows.Item(67).Insert()

$ws1.Range("A66").Value = "36A"
$ws1.Range("D66").Value = "36A"

$ws1.Range("A67").Value = "36B"
$ws1.Range("D67").Value = "36B"

# --- Keep the autofilter / used range in sync with the extra row --------
$ws1.AutoFilterMode = $false
$ws1.Range("A1:R78").AutoFilter()

# The workbook-level hidden "_xlnm._FilterDatabase" name for this sheet still
# points at the old range; move it out to the new last row too.
$filterName = $wb.Names.Item("1. Serotype!_FilterDatabase")
$filterName.RefersTo = "='1. Serotype'!`$A`$1:`$R`$78"

# --- View state: make "1. Serotype" the active/selected sheet -----------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("G81").Select()
